# Updates the Sema4b-Dcbld2 LR-pair sheet with refreshed TPM-derived
# expression / specificity / edge-weight numbers (commit: "update scripts
# wuth new tpm"). The sheet is a 5x5 cross of sending/target clusters
# (ECs, FAPs, Inflammatory-Mac, MuSCs, Resolving-Mac) for the single
# ligand-receptor pair Sema4b -> Dcbld2, rows 2..26 (5 sending blocks of
# 5 target rows each).
#
# Column layout:
#   G  Ligand average expression value           (per sending cluster)
#   H  Ligand total expression value              (per sending cluster)
#   I  Ligand derived specificity (avg)           = G / sum(G)
#   J  Ligand derived specificity (total)         = H / sum(H)
#   M  Receptor average expression value          (per target cluster)
#   N  Receptor total expression value            (per target cluster)
#   O  Receptor derived specificity (avg)         = M / sum(M)
#   P  Receptor derived specificity (total)       = N / sum(N)
#   Q  Edge average expression weight             = G * M
#   R  Edge total expression weight               = H * N
#   S  Edge average expression derived specificity= I * O
#   T  Edge total expression derived specificity  = J * P

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$clusters = @("ECs", "FAPs", "Inflammatory-Mac", "MuSCs", "Resolving-Mac")

$G = @{
    "ECs"               = 0.7257866666666666
    "FAPs"              = 2.719712
    "Inflammatory-Mac"  = 5.552337666666666
    "MuSCs"             = 4.039389
    "Resolving-Mac"     = 7.230285666666667
}
$H = @{
    "ECs"               = 2.17736
    "FAPs"              = 8.159136
    "Inflammatory-Mac"  = 16.657013
    "MuSCs"             = 8.078778
    "Resolving-Mac"     = 21.690857
}
$M = @{
    "ECs"               = 5.810536333333334
    "FAPs"              = 24.330847
    "Inflammatory-Mac"  = 0.5121463333333333
    "MuSCs"             = 10.7640175
    "Resolving-Mac"     = 0.6560336666666666
}
$N = @{
    "ECs"               = 17.431609
    "FAPs"              = 72.99254099999999
    "Inflammatory-Mac"  = 1.536439
    "MuSCs"             = 21.528035
    "Resolving-Mac"     = 1.968101
}

$sumG = 0.0
$sumH = 0.0
$sumM = 0.0
$sumN = 0.0
foreach ($c in $clusters) {
    $sumG += $G[$c]
    $sumH += $H[$c]
    $sumM += $M[$c]
    $sumN += $N[$c]
}

$row = 2
foreach ($send in $clusters) {
    $gVal = $G[$send]
    $hVal = $H[$send]
    $iVal = $gVal / $sumG
    $jVal = $hVal / $sumH

    foreach ($target in $clusters) {
        $mVal = $M[$target]
        $nVal = $N[$target]
        $oVal = $mVal / $sumM
        $pVal = $nVal / $sumN

        $qVal = $gVal * $mVal
        $rVal = $hVal * $nVal
        $sVal = $iVal * $oVal
        $tVal = $jVal * $pVal

        $ws.Cells.Item($row, 7).Value  = $gVal
        $ws.Cells.Item($row, 8).Value  = $hVal
        $ws.Cells.Item($row, 9).Value  = $iVal
        $ws.Cells.Item($row, 10).Value = $jVal
        $ws.Cells.Item($row, 13).Value = $mVal
        $ws.Cells.Item($row, 14).Value = $nVal
        $ws.Cells.Item($row, 15).Value = $oVal
        $ws.Cells.Item($row, 16).Value = $pVal
        $ws.Cells.Item($row, 17).Value = $qVal
        $ws.Cells.Item($row, 18).Value = $rVal
        $ws.Cells.Item($row, 19).Value = $sVal
        $ws.Cells.Item($row, 20).Value = $tVal

        $row++
    }
}
